$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.382.66"
$ws.Range("E2").Value = "  +4.34%  "
$ws.Range("D3").Value = "1.727.10"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'218.79"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'23.96"
$ws.Range("E8").Value = "  +3.18%  "
$ws.Range("D10").Value = "'0.0635"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "'0.0894"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.972.03"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "1.730.71"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").Value = "'4.23"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "'0.564"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "'67.55"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "28.324.12"
$ws.Range("E17").Value = "  +4.09%  "
$ws.Range("D18").Value = "'246.09"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "0.0₃0750"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'7.89"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "'9.64"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'149.25"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  +3.03%  "
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "1.485.54"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "'0.978"
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'0.601"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "'69.64"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'5.65"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "1.875.87"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").Value = "'1.72"
$ws.Range("E47").Value = "  +6.94%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'90.26"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0113"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "'8.16"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -1.02%  "
